# Added logging functionality: append new log rows to the worksheet,
# following the same "id in column A / timestamp in column D" layout
# already used by the existing log rows at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logEntries = @(
    @{ Row = 56; Id = "bf4d"; Timestamp = 45136.36827644676 },
    @{ Row = 57; Id = "3151"; Timestamp = 45136.36828064814 },
    @{ Row = 58; Id = "f699"; Timestamp = 45136.36829612269 }
)

# The last row already present in the sheet is used as the formatting
# template for each new row, so the appended cells reuse the exact same
# per-cell styling already used by the recent log entries: the id cell
# keeps the sheet's default (unstyled) look, and the timestamp cell keeps
# the custom date/time number format - instead of picking up the column
# default styles.
$templateRow = 55

foreach ($entry in $logEntries) {
    $r = $entry.Row

    # Ids are short hex-like tokens, but some (e.g. "3151") are made up of
    # digits only and would otherwise be auto-converted to a number by
    # COM's Value setter. Force text storage by pre-setting a text number
    # format, then restore the id cell's formatting from the template row
    # (applied uniformly for every id) so no stray number format/style is
    # ever left on the cell, and every id cell ends up unstyled, exactly
    # like the rest of the existing log rows.
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $entry.Id
    $ws.Range("A$templateRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Range("D$templateRow").Copy() | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($r, 4).Value = $entry.Timestamp

    $templateRow = $r
}

$excel.CutCopyMode = 0
